# Insert a new price record as row 525 in the "Repollo" sheet.
# Every existing row from 525 downward shifts down by one (Excel's
# native Rows.Insert behaviour), which reproduces the cascading diff
# where old row N's data now lives at row N+1 (and the former last
# "Segunda" row 649 becomes the new row 650).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 525..649 down to 526..650, carrying formatting from the
# row above (matches native Excel "insert" behaviour).
$ws.Rows.Item(525).Insert()

# Populate the newly-opened row 525 with the new record.
$ws.Range("A525").Value = 3
$ws.Range("B525").Value = "Femacal de La Calera"
$ws.Range("C525").Value = "Coquimbo"
$ws.Range("D525").Value = [DateTime]"2022-07-12"
$ws.Range("E525").Value = 5
$ws.Range("F525").Value = 100112006
$ws.Range("G525").Value = "Repollo"
$ws.Range("H525").Value = "Crespo record"
$ws.Range("I525").Value = "Primera"
$ws.Range("J525").Value = 3000
$ws.Range("K525").Value = 1300
$ws.Range("L525").Value = 1400
$ws.Range("M525").Value = 1360
$ws.Range("N525").Value = "$/unidad"
$ws.Range("O525").Value = "Provincia de Quillota"
$ws.Range("P525").Value = 1360
$ws.Range("Q525").Value = 1
$ws.Range("R525").Value = "Hortaliza"
